# "update some errors on data" -- corrects the recovered-count data on the
# "Recoverd" sheet and refreshes the active-sheet/selection view state that
# Excel records when a user is working on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recoverd")

# --- Data corrections (recovered counts were wrong) ---
$ws.Range("C21").Value = 0

$ws.Range("B28").Value = 26
$ws.Range("C28").Value = 1

$ws.Range("C29").Value = 4

$ws.Range("B30").Value = 33
$ws.Range("C30").Value = 3

# B31:B34 previously held running-total formulas (=SUM(prev+C)); replace
# them with the corrected plain values.
$ws.Range("B31").Value = 33
$ws.Range("C31").Value = 0

$ws.Range("B32").Value = 33

$ws.Range("B33").Value = 33

$ws.Range("B34").Value = 33

# --- View state: make "Recoverd" the active sheet/tab and update selection ---
$ws.Activate() | Out-Null
$ws.Range("C35").Select() | Out-Null
